$d = $word.ActiveDocument

# The heading run was originally split into "Re" + a collapsed _GoBack
# bookmark + "quest for repairs". Re-typing the whole heading merges
# those two runs back into a single "Request for repairs" run (and,
# since the bookmark fell inside the replaced range, removes it from
# here so it can be re-anchored further down, matching what Word does
# when the most recent edit point moves).
$d.Content.Find.Execute(
    "Request for repairs",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Request for repairs",
    2
)

# Only let the reader know they can reach the tenant "by phone, or by
# email" when a phone number actually exists -- wrap ", or by" in its
# own conditional instead of always showing it.
$d.Content.Find.Execute(
    ", or by email at",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{% if users[0].phone_numbers() %}, or by{% endif %} email at",
    2
)

# Word always tracks the location of the last edit with the hidden
# "_GoBack" bookmark. Re-anchor it at that last edit point (right after
# the {% endif %} that was just typed, before " email at").
$r = $d.Content
$r.Find.Execute("{% endif %} email at", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $r.Start + "{% endif %}".Length
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
